$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2646143.56
$ws.Range("C9").Value = 401026.99
$ws.Range("D9").Value = 3047170.55
$ws.Range("E9").Value = 13.16063487158604
$ws.Range("F9").Value = 86.83936512841397
$ws.Range("G9").Value = -61.24274102925761
$ws.Range("H9").Value = -52.21430491085437
$ws.Range("I9").Value = -53.63571654935976
$ws.Range("J9").Value = 26309
$ws.Range("K9").Value = 1120
$ws.Range("L9").Value = 27429
